$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, E and F hold values that must stay plain text (leading
# zeros / date-like timestamps / underscore+digits session ids), so
# force text format before assigning so Excel doesn't reinterpret them
# as numbers or dates, then clear the formatting again so the cell is
# left with the default (no explicit) style, matching a freshly
# authored row.
$textRange = $ws.Range("B8:B9,E8:E9,F8:F9")
$textRange.NumberFormat = "@"

$ws.Range("A8").Value = "subhan"
$ws.Range("B8").Value = "0987654321"
$ws.Range("C8").Value = "sd@sk.com"
$ws.Range("D8").Value = "pk"
$ws.Range("E8").Value = "2025-08-11 21:02:45"
$ws.Range("F8").Value = "session_1754928140"

$ws.Range("A9").Value = "subhan"
$ws.Range("B9").Value = "0987654432"
$ws.Range("C9").Value = "sad@sk.com"
$ws.Range("D9").Value = "wef"
$ws.Range("E9").Value = "2025-08-11 21:11:53"
$ws.Range("F9").Value = "session_1754928656"

$textRange.ClearFormats()
